# Scheduled-runner price/profit refresh for the Ixion_Profits workbook.
# Updates the currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H:N) for a handful of rows across each crafting-class sheet.

$wb = $excel.ActiveWorkbook

# --- ALC -------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H113").Value = 6669436.5
$ws.Range("I113").Value = 8335958
$ws.Range("J113").Value = 3350
$ws.Range("K113").Value = 8335958
$ws.Range("L113").Value = 3350
$ws.Range("M113").Value = -8332704
$ws.Range("N113").Value = -9858

$ws.Range("H116").Value = 18250.666
$ws.Range("I116").Value = 18250.666
$ws.Range("K116").Value = 18250.666
$ws.Range("M116").Value = -14808.666

$ws.Range("H129").Value = 990.2909
$ws.Range("J129").Value = 1095.9783
$ws.Range("L129").Value = 3287.9349
$ws.Range("N129").Value = -13287.9349

$ws.Range("H137").Value = 1547.5807
$ws.Range("I137").Value = 1338.3928
$ws.Range("J137").Value = 3500
$ws.Range("K137").Value = 4015.1784
$ws.Range("L137").Value = 10500
$ws.Range("M137").Value = -1465.1784
$ws.Range("N137").Value = -15600

$ws.Range("H138").Value = 1724.8636
$ws.Range("I138").Value = 997.5
$ws.Range("J138").Value = 2409.4412
$ws.Range("K138").Value = 2992.5
$ws.Range("L138").Value = 7228.323600000001
$ws.Range("M138").Value = 2147.5
$ws.Range("N138").Value = -17508.3236

# --- ARM -------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H101").Value = 70000
$ws.Range("J101").Value = 70000
$ws.Range("L101").Value = 70000
$ws.Range("N101").Value = -76490

$ws.Range("H122").Value = 1225403.2
$ws.Range("I122").Value = 1353718.6
$ws.Range("J122").Value = 6407
$ws.Range("K122").Value = 4061155.8
$ws.Range("L122").Value = 19221
$ws.Range("M122").Value = -4058705.8
$ws.Range("N122").Value = -24121

# --- BSM -------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H107").Value = 1132.0968
$ws.Range("I107").Value = 1067.375
$ws.Range("J107").Value = 1201.1333
$ws.Range("K107").Value = 1067.375
$ws.Range("L107").Value = 1201.1333
$ws.Range("M107").Value = 852.625
$ws.Range("N107").Value = -5041.1333

# --- CRP -------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H16").Value = 1769.6875
$ws.Range("I16").Value = 1500.6666
$ws.Range("K16").Value = 1500.6666
$ws.Range("M16").Value = -1213.6666

$ws.Range("H58").Value = 1599.2368
$ws.Range("I58").Value = 1286.4642
$ws.Range("J58").Value = 2475
$ws.Range("K58").Value = 1286.4642
$ws.Range("L58").Value = 2475
$ws.Range("M58").Value = -1083.4642
$ws.Range("N58").Value = -2881

$ws.Range("H94").Value = 4135.727
$ws.Range("I94").Value = 3796.1333
$ws.Range("J94").Value = 4418.722
$ws.Range("K94").Value = 3796.1333
$ws.Range("L94").Value = 4418.722
$ws.Range("M94").Value = -3345.1333
$ws.Range("N94").Value = -5320.722

$ws.Range("H99").Value = 5500.231
$ws.Range("J99").Value = 1500.5
$ws.Range("L99").Value = 1500.5
$ws.Range("N99").Value = -4496.5

$ws.Range("H113").Value = 1769.6875
$ws.Range("I113").Value = 1500.6666
$ws.Range("K113").Value = 1500.6666
$ws.Range("M113").Value = 669.3334

$ws.Range("H126").Value = 5500.231
$ws.Range("J126").Value = 1500.5
$ws.Range("L126").Value = 4501.5
$ws.Range("N126").Value = -9441.5

$ws.Range("H136").Value = 1599.2368
$ws.Range("I136").Value = 1286.4642
$ws.Range("J136").Value = 2475
$ws.Range("K136").Value = 3859.3926
$ws.Range("L136").Value = 7425
$ws.Range("M136").Value = -1309.3926
$ws.Range("N136").Value = -12525

# --- CUL -------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

$ws.Range("H40").Value = 435.27777
$ws.Range("I40").Value = 233.44827
$ws.Range("J40").Value = 1271.4286
$ws.Range("K40").Value = 933.79308
$ws.Range("L40").Value = 5085.7144
$ws.Range("M40").Value = -864.79308
$ws.Range("N40").Value = -5223.7144

$ws.Range("H86").Value = 1241
$ws.Range("I86").Value = 1300.9231
$ws.Range("J86").Value = 851.5
$ws.Range("K86").Value = 3902.7693
$ws.Range("L86").Value = 2554.5
$ws.Range("M86").Value = -2716.7693
$ws.Range("N86").Value = -4926.5

$ws.Range("H89").Value = 1241
$ws.Range("I89").Value = 1300.9231
$ws.Range("J89").Value = 851.5
$ws.Range("K89").Value = 11708.3079
$ws.Range("L89").Value = 7663.5
$ws.Range("M89").Value = -5780.3079
$ws.Range("N89").Value = -19519.5

# --- GSM -------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H33").Value = 5606.6665
$ws.Range("J33").Value = 5606.6665
$ws.Range("L33").Value = 5606.6665
$ws.Range("N33").Value = -6110.6665

$ws.Range("H40").Value = 12466.667
$ws.Range("J40").Value = 12466.667
$ws.Range("L40").Value = 12466.667
$ws.Range("N40").Value = -12768.667

$ws.Range("H102").Value = 998979.8
$ws.Range("I102").Value = 2120018
$ws.Range("J102").Value = 2501.4443
$ws.Range("K102").Value = 2120018
$ws.Range("L102").Value = 2501.4443
$ws.Range("M102").Value = -2118396
$ws.Range("N102").Value = -5745.4443

$ws.Range("H122").Value = 26622848
$ws.Range("J122").Value = 2905.9443
$ws.Range("L122").Value = 8717.832900000001
$ws.Range("N122").Value = -13617.8329

$ws.Range("H132").Value = 2787.18
$ws.Range("I132").Value = 3155.9583
$ws.Range("J132").Value = 2446.7693
$ws.Range("K132").Value = 9467.874899999999
$ws.Range("L132").Value = 7340.3079
$ws.Range("M132").Value = -6937.874899999999
$ws.Range("N132").Value = -12400.3079

# --- LTW -------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H40").Value = 83336270
$ws.Range("I40").Value = 83336270
$ws.Range("K40").Value = 83336270
$ws.Range("M40").Value = -83336134

$ws.Range("H122").Value = 5435472
$ws.Range("I122").Value = 5502083
$ws.Range("K122").Value = 16506249
$ws.Range("M122").Value = -16503799

$ws.Range("H132").Value = 8338347
$ws.Range("I132").Value = 9809202
$ws.Range("J132").Value = 3499.8333
$ws.Range("K132").Value = 29427606
$ws.Range("L132").Value = 10499.4999
$ws.Range("M132").Value = -29425076
$ws.Range("N132").Value = -15559.4999

# --- WVR -------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H3").Value = 26315.285
$ws.Range("J3").Value = 36801.4
$ws.Range("L3").Value = 36801.4
$ws.Range("N3").Value = -37029.4

$ws.Range("H81").Value = 1098.091
$ws.Range("I81").Value = 1107.9
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 2215.8
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -1154.8
$ws.Range("N81").Value = -4122

$ws.Range("H84").Value = 1098.091
$ws.Range("I84").Value = 1107.9
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 11079
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -5775
$ws.Range("N84").Value = -20608

$ws.Range("H122").Value = 1333.0714
$ws.Range("I122").Value = 1020.8
$ws.Range("J122").Value = 2113.75
$ws.Range("K122").Value = 3062.4
$ws.Range("L122").Value = 6341.25
$ws.Range("M122").Value = -612.3999999999996
$ws.Range("N122").Value = -11241.25

Write-Output "Applied scheduled profit refresh to ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR"
